# Update the "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets, matching the refreshed data pull from
# gh-pages output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 282
$ws1.Range("F4").Value  = 16601
$ws1.Range("F8").Value  = 359
$ws1.Range("F9").Value  = 207
$ws1.Range("F11").Value = 11550
$ws1.Range("F13").Value = 1231
$ws1.Range("F15").Value = 404
$ws1.Range("F17").Value = 62
$ws1.Range("F18").Value = 872
$ws1.Range("F20").Value = 148

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 282
$ws4.Range("F5").Value  = 16601
$ws4.Range("F9").Value  = 359
$ws4.Range("F10").Value = 207
$ws4.Range("F14").Value = 11550
$ws4.Range("F16").Value = 1231
$ws4.Range("F18").Value = 404
$ws4.Range("F20").Value = 62
$ws4.Range("F21").Value = 872
$ws4.Range("F23").Value = 148
